# Re-export of the playlist_per_work sheet: ids are regenerated from the
# speaker_variant text (no more reused/duplicated ids), rows are reordered
# to match the new export order, and the is_prefered ("x") flags are all
# cleared since the export no longer carries an is_pref / lev-distance
# selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row number -> (id, speaker_variant)
$data = @{
    2  = @("#de-geest", "De Geest")
    3  = @("#masinissa", "Masinissa")
    4  = @("#geest", "Geest")
    5  = @("#m.-lelivs", "M. Lelivs")
    6  = @("#sophonisba", "Sophonisba")
    7  = @("#hyppar", "Hyppar")
    8  = @("#amys", "Amys")
    9  = @("#masinis", "Masinis")
    10 = @("#sophon", "Sophon")
    11 = @("#priscvs", "Priscvs")
    12 = @("#scipi", "Scipi")
    13 = @("#gervnd", "Gervnd")
    14 = @("#scipio", "Scipio")
    15 = @("#masi", "Masi")
    16 = @("#amystas", "Amystas")
    17 = @("#masin", "Masin")
    18 = @("#amyst", "Amyst")
    19 = @("#gervndvla", "Gervndvla")
    20 = @("#rey", "REY")
    21 = @("#iongen", "Iongen")
    22 = @("#lelivs", "Lelivs")
    23 = @("#de-geest-sophonisba", "De Geest Sophonisba")
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = ""
}
